$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V between row 19 and row 20 (teams/odds reordered) ---
$ws.Range("F19").Value2 = "Vozdovac"
$ws.Range("F20").Value2 = "TSC"
$ws.Range("G19").Value2 = 1
$ws.Range("G20").Value2 = 1
$ws.Range("H19").Value2 = "Radnik"
$ws.Range("H20").Value2 = "Radnicki 1923"
$ws.Range("I19").Value2 = 1
$ws.Range("I20").Value2 = 0
$ws.Range("J19").Value2 = 3.21
$ws.Range("J20").Value2 = 1.33
$ws.Range("K19").Value2 = "11/08/2023 07:12"
$ws.Range("K20").Value2 = "11/08/2023 07:12"
$ws.Range("L19").Value2 = 2.12
$ws.Range("L20").Value2 = 1.47
$ws.Range("M19").Value2 = "12/08/2023 18:54"
$ws.Range("M20").Value2 = "12/08/2023 16:58"
$ws.Range("N19").Value2 = 3.12
$ws.Range("N20").Value2 = 4.58
$ws.Range("O19").Value2 = "11/08/2023 07:12"
$ws.Range("O20").Value2 = "11/08/2023 07:12"
$ws.Range("P19").Value2 = 3.16
$ws.Range("P20").Value2 = 4.35
$ws.Range("Q19").Value2 = "12/08/2023 18:54"
$ws.Range("Q20").Value2 = "12/08/2023 18:50"
$ws.Range("R19").Value2 = 2.11
$ws.Range("R20").Value2 = 7
$ws.Range("S19").Value2 = "11/08/2023 07:12"
$ws.Range("S20").Value2 = "11/08/2023 07:12"
$ws.Range("T19").Value2 = 3.53
$ws.Range("T20").Value2 = 6.26
$ws.Range("U19").Value2 = "12/08/2023 18:54"
$ws.Range("U20").Value2 = "12/08/2023 18:50"
$ws.Range("V19").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-radnik-surdulica/WOOL9vKS/"
$ws.Range("V20").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-radnicki-1923/vVDe3xcd/"

# --- Swap F:V between row 38 and row 39 (teams/odds reordered) ---
$ws.Range("F38").Value2 = "Sp. Subotica"
$ws.Range("F39").Value2 = "Cukaricki"
$ws.Range("G38").Value2 = 1
$ws.Range("G39").Value2 = 0
$ws.Range("H38").Value2 = "Napredak"
$ws.Range("H39").Value2 = "Radnik"
$ws.Range("I38").Value2 = 3
$ws.Range("I39").Value2 = 0
$ws.Range("J38").Value2 = 1.93
$ws.Range("J39").Value2 = 1.44
$ws.Range("K38").Value2 = "24/08/2023 09:13"
$ws.Range("K39").Value2 = "24/08/2023 09:13"
$ws.Range("L38").Value2 = 2
$ws.Range("L39").Value2 = 1.4
$ws.Range("M38").Value2 = "27/08/2023 19:29"
$ws.Range("M39").Value2 = "27/08/2023 19:25"
$ws.Range("N38").Value2 = 3.1
$ws.Range("N39").Value2 = 4.05
$ws.Range("O38").Value2 = "24/08/2023 09:13"
$ws.Range("O39").Value2 = "24/08/2023 09:13"
$ws.Range("P38").Value2 = 3.23
$ws.Range("P39").Value2 = 4.26
$ws.Range("Q38").Value2 = "27/08/2023 19:29"
$ws.Range("Q39").Value2 = "27/08/2023 19:25"
$ws.Range("R38").Value2 = 3.63
$ws.Range("R39").Value2 = 5.74
$ws.Range("S38").Value2 = "24/08/2023 09:13"
$ws.Range("S39").Value2 = "24/08/2023 09:13"
$ws.Range("T38").Value2 = 3.82
$ws.Range("T39").Value2 = 8.359999999999999
$ws.Range("U38").Value2 = "27/08/2023 19:29"
$ws.Range("U39").Value2 = "27/08/2023 19:25"
$ws.Range("V38").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/spartak-subotica-napredak/rNGLvaJ2/"
$ws.Range("V39").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-radnik-surdulica/Iu9Yyc3L/"

# --- Swap F:V between row 50 and row 51 (teams/odds reordered) ---
$ws.Range("F50").Value2 = "Vozdovac"
$ws.Range("F51").Value2 = "TSC"
$ws.Range("G50").Value2 = 3
$ws.Range("G51").Value2 = 6
$ws.Range("H50").Value2 = "Mladost"
$ws.Range("H51").Value2 = "Zeleznicar Pancevo"
$ws.Range("I50").Value2 = 1
$ws.Range("I51").Value2 = 3
$ws.Range("J50").Value2 = 2.07
$ws.Range("J51").Value2 = 1.21
$ws.Range("K50").Value2 = "14/09/2023 09:13"
$ws.Range("K51").Value2 = "14/09/2023 09:13"
$ws.Range("L50").Value2 = 1.83
$ws.Range("L51").Value2 = 1.33
$ws.Range("M50").Value2 = "16/09/2023 18:52"
$ws.Range("M51").Value2 = "16/09/2023 18:52"
$ws.Range("N50").Value2 = 3.13
$ws.Range("N51").Value2 = 5.61
$ws.Range("O50").Value2 = "14/09/2023 09:13"
$ws.Range("O51").Value2 = "14/09/2023 09:13"
$ws.Range("P50").Value2 = 3.58
$ws.Range("P51").Value2 = 4.42
$ws.Range("Q50").Value2 = "16/09/2023 18:52"
$ws.Range("Q51").Value2 = "16/09/2023 18:52"
$ws.Range("R50").Value2 = 3.2
$ws.Range("R51").Value2 = 9.23
$ws.Range("S50").Value2 = "14/09/2023 09:13"
$ws.Range("S51").Value2 = "14/09/2023 09:13"
$ws.Range("T50").Value2 = 4.06
$ws.Range("T51").Value2 = 10.82
$ws.Range("U50").Value2 = "16/09/2023 18:52"
$ws.Range("U51").Value2 = "16/09/2023 18:52"
$ws.Range("V50").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-mladost-lucani/25QqVon6/"
$ws.Range("V51").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-zeleznicar-pancevo/xOIdSqWO/"

# --- Swap F:V between row 84 and row 85 (teams/odds reordered) ---
$ws.Range("F84").Value2 = "Vojvodina"
$ws.Range("F85").Value2 = "Mladost"
$ws.Range("G84").Value2 = 3
$ws.Range("G85").Value2 = 1
$ws.Range("H84").Value2 = "Radnik"
$ws.Range("H85").Value2 = "Radnicki Nis"
$ws.Range("I84").Value2 = 0
$ws.Range("I85").Value2 = 2
$ws.Range("J84").Value2 = 1.61
$ws.Range("J85").Value2 = 1.78
$ws.Range("K84").Value2 = "20/10/2023 01:42"
$ws.Range("K85").Value2 = "20/10/2023 01:42"
$ws.Range("L84").Value2 = 1.34
$ws.Range("L85").Value2 = 2.87
$ws.Range("M84").Value2 = "22/10/2023 17:51"
$ws.Range("M85").Value2 = "22/10/2023 17:59"
$ws.Range("N84").Value2 = 3.58
$ws.Range("N85").Value2 = 3.36
$ws.Range("O84").Value2 = "20/10/2023 01:42"
$ws.Range("O85").Value2 = "20/10/2023 01:42"
$ws.Range("P84").Value2 = 4.63
$ws.Range("P85").Value2 = 3.13
$ws.Range("Q84").Value2 = "22/10/2023 17:51"
$ws.Range("Q85").Value2 = "22/10/2023 17:59"
$ws.Range("R84").Value2 = 4.69
$ws.Range("R85").Value2 = 4.07
$ws.Range("S84").Value2 = "20/10/2023 01:42"
$ws.Range("S85").Value2 = "20/10/2023 01:42"
$ws.Range("T84").Value2 = 9.609999999999999
$ws.Range("T85").Value2 = 2.48
$ws.Range("U84").Value2 = "22/10/2023 17:51"
$ws.Range("U85").Value2 = "22/10/2023 17:59"
$ws.Range("V84").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/vojvodina-radnik-surdulica/Ctudti6r/"
$ws.Range("V85").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/mladost-lucani-radnicki-nis/0GjCx957/"

# --- Swap F:V between row 90 and row 91 (teams/odds reordered) ---
$ws.Range("F90").Value2 = "IMT Novi Beograd"
$ws.Range("F91").Value2 = "Radnicki Nis"
$ws.Range("G90").Value2 = 1
$ws.Range("G91").Value2 = 1
$ws.Range("H90").Value2 = "Crvena zvezda"
$ws.Range("H91").Value2 = "Sp. Subotica"
$ws.Range("I90").Value2 = 2
$ws.Range("I91").Value2 = 1
$ws.Range("J90").Value2 = 8.15
$ws.Range("J91").Value2 = 1.7
$ws.Range("K90").Value2 = "27/10/2023 06:42"
$ws.Range("K91").Value2 = "27/10/2023 06:42"
$ws.Range("L90").Value2 = 24.2
$ws.Range("L91").Value2 = 1.66
$ws.Range("M90").Value2 = "28/10/2023 18:29"
$ws.Range("M91").Value2 = "28/10/2023 18:23"
$ws.Range("N90").Value2 = 5.6
$ws.Range("N91").Value2 = 3.42
$ws.Range("O90").Value2 = "27/10/2023 06:42"
$ws.Range("O91").Value2 = "27/10/2023 06:42"
$ws.Range("P90").Value2 = 9.529999999999999
$ws.Range("P91").Value2 = 3.65
$ws.Range("Q90").Value2 = "28/10/2023 18:29"
$ws.Range("Q91").Value2 = "28/10/2023 18:23"
$ws.Range("R90").Value2 = 1.23
$ws.Range("R91").Value2 = 4.23
$ws.Range("S90").Value2 = "27/10/2023 06:42"
$ws.Range("S91").Value2 = "27/10/2023 06:42"
$ws.Range("T90").Value2 = 1.09
$ws.Range("T91").Value2 = 5.08
$ws.Range("U90").Value2 = "28/10/2023 18:21"
$ws.Range("U91").Value2 = "28/10/2023 18:23"
$ws.Range("V90").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-crvena-zvezda/SjAgknkD/"
$ws.Range("V91").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-spartak-subotica/2qDshl5f/"

# --- Swap F:V between row 122 and row 123 (teams/odds reordered) ---
$ws.Range("F122").Value2 = "Mladost"
$ws.Range("F123").Value2 = "Partizan"
$ws.Range("G122").Value2 = 1
$ws.Range("G123").Value2 = 3
$ws.Range("H122").Value2 = "Zeleznicar Pancevo"
$ws.Range("H123").Value2 = "Vojvodina"
$ws.Range("I122").Value2 = 0
$ws.Range("I123").Value2 = 1
$ws.Range("J122").Value2 = 2.32
$ws.Range("J123").Value2 = 1.5
$ws.Range("K122").Value2 = "01/12/2023 03:43"
$ws.Range("K123").Value2 = "01/12/2023 03:43"
$ws.Range("L122").Value2 = 2.14
$ws.Range("L123").Value2 = 1.51
$ws.Range("M122").Value2 = "02/12/2023 15:03"
$ws.Range("M123").Value2 = "02/12/2023 15:26"
$ws.Range("N122").Value2 = 2.97
$ws.Range("N123").Value2 = 4.01
$ws.Range("O122").Value2 = "01/12/2023 03:43"
$ws.Range("O123").Value2 = "01/12/2023 03:43"
$ws.Range("P122").Value2 = 3.24
$ws.Range("P123").Value2 = 3.71
$ws.Range("Q122").Value2 = "02/12/2023 15:03"
$ws.Range("Q123").Value2 = "02/12/2023 15:26"
$ws.Range("R122").Value2 = 2.94
$ws.Range("R123").Value2 = 5.3
$ws.Range("S122").Value2 = "01/12/2023 03:43"
$ws.Range("S123").Value2 = "01/12/2023 03:43"
$ws.Range("T122").Value2 = 3.38
$ws.Range("T123").Value2 = 7.24
$ws.Range("U122").Value2 = "02/12/2023 15:03"
$ws.Range("U123").Value2 = "02/12/2023 15:26"
$ws.Range("V122").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/mladost-lucani-zeleznicar-pancevo/zka4Yd9c/"
$ws.Range("V123").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/partizan-vojvodina/tUjxVN25/"

# --- Append new row 144 (new match: Radnicki 1923 vs Radnik) ---
$ws.Range("A143:V143").Copy($ws.Range("A144:V144"))
$ws.Range("A144").Value2 = 143
$ws.Range("B144").Value2 = "serbia"
$ws.Range("C144").Value2 = "super-liga"
$ws.Range("D144").Value2 = "2023-2024"
$ws.Range("E144").Value2 = 45280.54166666666
$ws.Range("F144").Value2 = "Radnicki 1923"
$ws.Range("G144").Value2 = 1
$ws.Range("H144").Value2 = "Radnik"
$ws.Range("I144").Value2 = 0
$ws.Range("J144").Value2 = 2.04
$ws.Range("K144").Value2 = "25/09/2023 03:12"
$ws.Range("L144").Value2 = 1.55
$ws.Range("M144").Value2 = "20/12/2023 12:53"
$ws.Range("N144").Value2 = 2.98
$ws.Range("O144").Value2 = "25/09/2023 03:12"
$ws.Range("P144").Value2 = 3.66
$ws.Range("Q144").Value2 = "20/12/2023 12:53"
$ws.Range("R144").Value2 = 3.45
$ws.Range("S144").Value2 = "25/09/2023 03:12"
$ws.Range("T144").Value2 = 6.57
$ws.Range("U144").Value2 = "20/12/2023 12:53"
$ws.Range("V144").Value2 = "https://www.betexplorer.com/football/serbia/super-liga/radnicki-1923-radnik-surdulica/4CHCgMUN/"
